$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 533 (shifts old rows 533-619 down to
# 534-620, inheriting formatting -- including the date number-format on
# column D -- from the row above, same as Excel's native "Insert Row").
$ws.Rows.Item(533).Insert()

# Populate the newly inserted row 533 with the new weekly data point.
# J533 and L533 keep the same values the (now shifted-down) old row 533
# had, so we simply copy them forward from row 534 (which holds what used
# to be row 533's data after the shift).
$ws.Range("A533").Value = $ws.Range("A534").Value()
$ws.Range("B533").Value = $ws.Range("B534").Value()
$ws.Range("C533").Value = $ws.Range("C534").Value()
$ws.Range("D533").Value = 45218
$ws.Range("E533").Value = $ws.Range("E534").Value()
$ws.Range("F533").Value = $ws.Range("F534").Value()
$ws.Range("G533").Value = $ws.Range("G534").Value()
$ws.Range("H533").Value = $ws.Range("H534").Value()
$ws.Range("I533").Value = $ws.Range("I534").Value()
$ws.Range("J533").Value = $ws.Range("J534").Value()
$ws.Range("K533").Value = 5800
$ws.Range("L533").Value = $ws.Range("L534").Value()
$ws.Range("M533").Value = 5900
$ws.Range("N533").Value = $ws.Range("N534").Value()
$ws.Range("O533").Value = $ws.Range("O534").Value()
$ws.Range("P533").Value = 295
$ws.Range("Q533").Value = $ws.Range("Q534").Value()
$ws.Range("R533").Value = $ws.Range("R534").Value()
